$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.221.46"
$ws.Range("E2").Value = '  -3.30%  '
$ws.Range("D3").Value = "'2.994.39"
$ws.Range("E3").Value = '  -4.23%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'580.55"
$ws.Range("E5").Value = '  -2.47%  '
$ws.Range("D6").Value = "'145.94"
$ws.Range("E6").Value = '  -8.27%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'0.522"
$ws.Range("E8").Value = '  -3.80%  '
$ws.Range("D9").Value = "'2.997.59"
$ws.Range("E9").Value = '  -4.03%  '
$ws.Range("E10").Value = '  -7.41%  '
$ws.Range("E11").Value = '  -5.46%  '
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = '  -2.94%  '
$ws.Range("D13").Value = "'0.0000227"
$ws.Range("E13").Value = '  -5.98%  '
$ws.Range("D14").Value = "'34.53"
$ws.Range("E14").Value = '  -7.79%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = "'3.493.84"
$ws.Range("E16").Value = '  -3.92%  '
$ws.Range("D17").Value = "'7.06"
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").Value = "'62.261.35"
$ws.Range("E18").Value = '  -3.02%  '
$ws.Range("D19").Value = "'2.997.48"
$ws.Range("E19").Value = '  -3.94%  '
$ws.Range("D20").Value = "'456.76"
$ws.Range("E20").Value = '  -4.87%  '
$ws.Range("D21").Value = "'13.85"
$ws.Range("E21").Value = '  -5.22%  '
$ws.Range("D22").Value = "'0.677"
$ws.Range("E22").Value = '  -5.69%  '
$ws.Range("D23").Value = "'7.28"
$ws.Range("E23").Value = '  -4.80%  '
$ws.Range("D24").Value = "'80.03"
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("E25").Value = '  -8.51%  '
$ws.Range("D26").Value = "'12.24"
$ws.Range("E26").Value = '  -6.39%  '
$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = '  -5.88%  '
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = '  -5.11%  '
$ws.Range("D31").Value = "'2.60"
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("E32").Value = '  -6.26%  '
$ws.Range("D33").Value = "'26.90"
$ws.Range("E34").Value = '  -5.86%  '
$ws.Range("E35").Value = '  -4.48%  '
$ws.Range("D36").Value = "'0.0₃0781"
$ws.Range("E36").Value = '  -8.44%  '
$ws.Range("D37").Value = "'5.72"
$ws.Range("E37").Value = '  -5.70%  '
$ws.Range("E38").Value = '  -6.75%  '
$ws.Range("D39").Value = "'50.02"
$ws.Range("E39").Value = '  -2.19%  '
$ws.Range("D40").Value = "'9.00"
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = '  -12.89%  '
$ws.Range("D42").Value = "'410.03"
$ws.Range("E42").Value = '  -9.73%  '
$ws.Range("E43").Value = '  -6.12%  '
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = "'2.766.98"
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("D46").Value = "'0.0350"
$ws.Range("E46").Value = '  -4.74%  '
$ws.Range("D47").Value = "'38.51"
$ws.Range("E47").Value = '  -3.97%  '
$ws.Range("D48").Value = "'127.60"
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").Value = '  -2.83%  '
$ws.Range("D51").Value = "'23.70"
$ws.Range("E51").Value = '  -8.57%  '
